# Adds a short lead-in sentence (with an italic date) directly above the
# chart image in the "Hospitalisation" and "ICU" sections, and a short
# plain-text lead-in above the chart image in the "Number of confirmed
# infections" section. In each case the image paragraph's style changes
# from "First Paragraph" to "Body Text", matching the pattern already
# used elsewhere in the document (e.g. the "How trustworthy is this?"
# section).

function Insert-DatedParagraph($d, $imgParaIndex, $prefixText, $dateText, $suffixText) {
    $imgPara = $d.Paragraphs.Item($imgParaIndex)
    $imgPara.Range.InsertParagraphBefore()

    # The freshly inserted (empty) paragraph now occupies $imgParaIndex;
    # the original image paragraph has shifted to $imgParaIndex + 1.
    $newPara = $d.Paragraphs.Item($imgParaIndex)
    $newRange = $newPara.Range
    $newRange.Collapse(1)          # wdCollapseStart
    $newRange.InsertAfter($prefixText)
    $newRange.Collapse(0)          # wdCollapseEnd
    $newRange.InsertAfter(" ")
    $newRange.Collapse(0)
    $dateStart = $newRange.Start
    $newRange.InsertAfter($dateText)
    $dateEnd = $newRange.End
    $newRange.Collapse(0)
    $newRange.InsertAfter($suffixText)

    $dateRange = $d.Range($dateStart, $dateEnd)
    $dateRange.Font.Italic = 1

    $imgParaAfter = $d.Paragraphs.Item($imgParaIndex + 1)
    $imgParaAfter.Style = "Body Text"
}

function Insert-PlainParagraph($d, $imgParaIndex, $text) {
    $imgPara = $d.Paragraphs.Item($imgParaIndex)
    $imgPara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($imgParaIndex)
    $newRange = $newPara.Range
    $newRange.Collapse(1)          # wdCollapseStart
    $newRange.InsertAfter($text)

    $imgParaAfter = $d.Paragraphs.Item($imgParaIndex + 1)
    $imgParaAfter.Style = "Body Text"
}

$d = $word.ActiveDocument

# Process from the bottom of the document upward so that earlier
# (not-yet-processed) paragraph indices are unaffected by the
# paragraphs inserted later in the script.

# "Number of confirmed infections" section -> image paragraph #13
Insert-PlainParagraph $d 13 "It is not possible to predict accurately when the current outbreak will peak. It is too far in the future."

# "ICU" section -> image paragraph #9
Insert-DatedParagraph $d 9 "Every ICU bed will be occupied on on" "2022-02-14" "."

# "Hospitalisation" section -> image paragraph #7
Insert-DatedParagraph $d 7 "Hospitals will be saturated on" "2022-02-17" "."
